$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O header: "Internal Assignment" (bold, size 12 - like the other
# bold sub-header cells K4:N4, just one point bigger)
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Size = 12

# New column O data cells: literal text "TRUE"/"FALSE" (not boolean values),
# using the same plain style as the rest of the data rows (column A etc).
# A scratch cell (M5, currently empty) is used to get Excel to store the
# value as literal text (quote-prefix) instead of auto-converting it to a
# boolean, the value is then copied over and the scratch cell cleared.
$ws.Range("M5").Value = "'TRUE"
$ws.Range("M5").Copy()
$ws.Range("O5").PasteSpecial(-4163)
$ws.Range("M5").Clear()

$ws.Range("M5").Value = "'FALSE"
$ws.Range("M5").Copy()
$ws.Range("O6").PasteSpecial(-4163)
$ws.Range("O7").PasteSpecial(-4163)
$ws.Range("O8").PasteSpecial(-4163)
$ws.Range("M5").Clear()

$ws.Range("A5").Copy()
$ws.Range("O5:O8").PasteSpecial(-4122)

$ws.Range("A1").Select()
